$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# repull data, push all data, mean calculation
# Update dSF (F column) values to match re-pulled data
$ws.Range("F2").Value = -7
$ws.Range("F6").Value = 1
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = 0
